$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 367, shifting existing rows 367-416 down to 368-417
$ws.Rows.Item(367).Insert()

# Populate the newly inserted row 367 with the new weekly record
$ws.Range("A367").Value = 4
$ws.Range("B367").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C367").Value = "Los Lagos"
$ws.Range("D367").Value = 45131
$ws.Range("E367").Value = 10
$ws.Range("F367").Value = "Fruta"
$ws.Range("G367").Value = 100108
$ws.Range("H367").Value = "Tropicales y subtropicales"
$ws.Range("I367").Value = 100108002
$ws.Range("J367").Value = "Mango"
$ws.Range("K367").Value = "Sin especificar"
$ws.Range("L367").Value = "Primera"
$ws.Range("M367").Value = 200
$ws.Range("N367").Value = 9500
$ws.Range("O367").Value = 10000
$ws.Range("P367").Value = 9750
$ws.Range("Q367").Value = "$/bandeja 4 kilos"
$ws.Range("R367").Value = "Perú"
$ws.Range("S367").Value = 2438
$ws.Range("T367").Value = 4
